$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Insert a new column before column C ("warehouse") ---------------------
# This shifts every existing column C..AI one slot right (-> D..AJ), carrying
# along each row's style, which is exactly what the target diff shows.
$ws.Columns("C:C").Insert()

# --- Populate the new "warehouse" column ------------------------------------
$ws.Range("C3").Value = "warehouse"
$ws.Range("C4").Value = "whrm01"
$ws.Range("C5").Value = "whrm02"

# --- Column widths -----------------------------------------------------------
# Column B ("partno") grew slightly, and the new warehouse column C needs its
# own explicit width.
$ws.Columns("B:B").ColumnWidth = 27.5
$ws.Columns("C:C").ColumnWidth = 22.43

# --- Selection / view ---------------------------------------------------------
# Author scrolled back to the left (no more frozen/left-anchored topLeftCell)
# and left the selection on D18 instead of the old J5:AI5 block.
$ws.Range("D18").Select()
